# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet, which contain the same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8809
    $ws.Range("F3").Value = 199
    $ws.Range("F4").Value = 427
    $ws.Range("F5").Value = 341
}
